$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert a new row at position 12 (pushes the existing rows 12-46 down to 13-47,
# carrying all of their values/styles/merges along automatically)
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new "compensacion administrativa
# retroactiva" catalog entry
$ws.Cells.Item(12, 1).Value = "COMPENSACION_ADMINISTRATIVA_RETRO"
$ws.Cells.Item(12, 2).Value = 1
$ws.Cells.Item(12, 3).Value = 1102
$ws.Cells.Item(12, 4).Value = "COMPENSACION ADMINISTRATIVA RETROACTIVA"
$ws.Cells.Item(12, 5).Value = "'038"
$ws.Cells.Item(12, 6).Value = 0

# Reflect the scrolled / selected view state from the edit
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Application.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D12").Select()
